$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet by copying the existing
#    "2022-Q3" sheet (same column layout/styles), placed right
#    before it, then renaming and overwriting its data.
# ------------------------------------------------------------------
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$srcQ3.Copy($srcQ3)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

$fundRows = @(
    @(0, '011329', '景顺长城新能源产业股票C', '43.23', '90.68', '4.24', '1.8330', 9),
    @(1, '260101', '景顺长城优选混合', '44.73', '76.26', '2.63', '1.1764', 8),
    @(2, '011328', '景顺长城新能源产业股票A', '24.53', '90.68', '4.24', '1.0401', 9),
    @(3, '213003', '宝盈策略增长混合', '10.36', '91.66', '7.99', '0.8278', 4),
    @(4, '213002', '宝盈泛沿海增长混合', '5.08', '92.58', '9.31', '0.4729', 2),
    @(5, '009376', '景顺长城成长领航混合', '12.23', '92.99', '3.27', '0.3999', 9),
    @(6, '260111', '景顺长城公司治理混合', '3.58', '91.54', '3.28', '0.1174', 9),
    @(7, '000796', '宝盈睿丰创新灵活配置混合 - C', '0.61', '92.26', '9.14', '0.0558', 2),
    @(8, '000794', '宝盈睿丰创新灵活配置混合 - A/B', '0.39', '92.26', '9.14', '0.0356', 2),
    @(9, '006644', '弘毅远方消费升级混合A', '0.39', '84.03', '2.37', '0.0092', 6),
    @(10, '011438', '红塔红土盛昌优选混合A', '0.13', '92.67', '4.85', '0.0063', 5),
    @(11, '001535', '景顺长城改革机遇灵活配置混合A', '0.28', '43.79', '1.62', '0.0045', 8),
    @(12, '007945', '景顺长城改革机遇灵活配置混合C', '0.24', '43.79', '1.62', '0.0039', 8),
    @(13, '011439', '红塔红土盛昌优选混合C', '0.04', '92.67', '4.85', '0.0019', 5),
    @(14, '014422', '弘毅远方消费升级混合C', '0.06', '84.03', '2.37', '0.0014', 6),
    @(15, '004725', '先锋聚元灵活配置混合C', '0.04', '92.62', '2.67', '0.0011', 2),
    @(16, '004724', '先锋聚元灵活配置混合A', '0.03', '92.62', '2.67', '0.0008', 2),
    @(17, '003587', '先锋精一灵活配置混合C', '0.02', '93.29', '3.19', '0.0006', 1),
    @(18, '004727', '先锋聚优灵活配置混合C', '0.02', '93.09', '2.56', '0.0005', 4),
    @(19, '004726', '先锋聚优灵活配置混合A', '0.01', '93.09', '2.56', '0.0003', 4),
    @(20, '003586', '先锋精一灵活配置混合A', '0.01', '93.29', '3.19', '0.0003', 1),

)

# Remove the extra rows (2022-Q3 has 28 data rows, 2022-Q4 has 22 data
# rows -> rows 23-28 are dropped so the sheet ends up A1:H22).
$newSheet.Rows("23:28").Delete()

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $row[0]

    $newSheet.Range("B$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[1]

    $newSheet.Range("C$r").NumberFormat = "@"
    $newSheet.Range("C$r").Value = $row[2]

    $newSheet.Range("D$r").NumberFormat = "@"
    $newSheet.Range("D$r").Value = $row[3]

    $newSheet.Range("E$r").NumberFormat = "@"
    $newSheet.Range("E$r").Value = $row[4]

    $newSheet.Range("F$r").NumberFormat = "@"
    $newSheet.Range("F$r").Value = $row[5]

    $newSheet.Range("G$r").NumberFormat = "@"
    $newSheet.Range("G$r").Value = $row[6]

    $newSheet.Range("H$r").Value = $row[7]

    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row 2 holding the
#    2022-Q4 totals, push the rest down, and renumber the index
#    column (A).
# ------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()

# Re-apply formatting to the freshly inserted row by copying it from
# the row below (which still holds the old row-2 formatting/values
# that got pushed down to row 3).
$totalWs.Range("A3:D3").Copy($totalWs.Range("A2:D2"))

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q4"
$totalWs.Range("C2").Value = 21
$totalWs.Range("D2").Value = 5.99

# Renumber the index column A for the rows that shifted down.
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
$totalWs.Range("A7").Value = 5

Write-Host "done"
